$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting of the existing station rows (row2 = station A, row3 = station B)
# into the two new rows so the same cell styles (date/number formats) are reused.
$ws.Range("A2:F2").Copy()
$ws.Range("A18:F18").PasteSpecial(-4122)

$ws.Range("A3:F3").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)

# Row 18: 四方坪站 (station A) for 2026-01-09 (serial 46031)
$ws.Range("A18").Value = 46031
$ws.Range("B18").Value = "四方坪站"
$ws.Range("C18").Value = 15021.67
$ws.Range("D18").Value = 10154.33
$ws.Range("E18").Value = 3394.87
$ws.Range("F18").Value = 642

# Row 19: 高岭站 (station B) for 2026-01-09 (serial 46031)
$ws.Range("A19").Value = 46031
$ws.Range("B19").Value = "高岭站"
$ws.Range("C19").Value = 5035.63
$ws.Range("D19").Value = 4126.6899999999996
$ws.Range("E19").Value = 1375.28
$ws.Range("F19").Value = 179

# Update selection to match the diff
$ws.Range("H20").Select()
